# Apply the BOM update described in the commit:
# "added pul-ups for encoder A,B,Switch, board rerouted, euricircuits pass and new auto quotation"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PurchaseList")

# --- Designator text updates -------------------------------------------------
# Header row designators renamed (JP1,JP2 -> J2,J3 ; JP3 -> J1 ; LCD1 -> J4)
$ws.Range("B4").Value = "J2, J3"
$ws.Range("B5").Value = "R1, R3, R4, R5"
$ws.Range("B9").Value = "J1"
$ws.Range("B10").Value = "J4"

# Row 9 used to be the header connector (HEADER_5 / 1X05) and is now the LCD
# module row; row 10 used to be the LCD module and is now a plain header.
$ws.Range("C9").Value = "LCD MODULE SIL"
$ws.Range("D9").Value = "1X14-FEMALE"
$ws.Range("E9").Value = "LCD MODULE SINGLE IN LINE"

$ws.Range("C10").Value = "HEADER_5"
$ws.Range("D10").Value = "1X05"
$ws.Range("E10").Value = ""

# --- Report time stamp --------------------------------------------------------
$ws.Range("F14").Value = "21:45"

# --- New auto quotation (supplier stock / unit prices) -----------------------
$ws.Range("J2").Value = 101930
$ws.Range("P2").Value = 0.09016

$ws.Range("J4").Value = 2985
$ws.Range("P4").Value = 0.49679

$ws.Range("J5").Value = 6950
$ws.Range("M5").Value = 4
$ws.Range("N5").Value = 20
$ws.Range("P5").Value = 0.09016

$ws.Range("P7").Value = 0.67621

# Recalculate so dependent formulas (Q column, G.Total, per-board, NOW()) update
$excel.CalculateFullRebuild()
$wb.Save()
